# Converts an EMU length into the "points" value to feed into the
# PowerPoint object model (Shape.Width / Shape.Height are in points,
# 1 pt == 12700 EMU). A small epsilon (half an EMU, in points) is
# added so the point -> EMU re-conversion performed internally lands
# squarely on the target EMU value instead of the adjacent one below it.
function EmuToPt($emu) {
    return ($emu / 12700) + (0.5 / 12700)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cr = [char]13
$dash = [char]0x2013

# ---------------------------------------------------------------
# Shape 27 ("TextBox 75") - the "3. clone" box.
#   * "git clone https://" -> "git clone git@github.com:"
#   * drop the literal "@github.com/" + the 2nd italic "user" run
#   * "git remote add upstream https://github.com/ibm/lale.git"
#       -> "git remote add upstream git@github.com:ibm/lale.git"
#   * shrink the box width to fit the (now shorter) text
# ---------------------------------------------------------------
$cloneShape = $s.Shapes.Item(27)
$cloneRange = $cloneShape.TextFrame.TextRange

# Paragraph 4: "git remote add upstream https://github.com/ibm/lale.git"
# (55 chars, starting at char 86) becomes a single run with the new text.
$cloneRange.Characters(86, 55).Text = "git remote add upstream git@github.com:ibm/lale.git"

# Paragraph 2: delete the "@github.com/" run + the trailing italic "user"
# run (16 chars total, starting at char 32) so "user"/"lale.git" abut.
$cloneRange.Characters(32, 16).Text = ""

# Paragraph 2: "git clone https://" (18 chars, starting at char 10)
# becomes "git clone git@github.com:".
$cloneRange.Characters(10, 18).Text = "git clone git@github.com:"

$cloneShape.Width = EmuToPt 2571538

# ---------------------------------------------------------------
# Shape 35 ("TextBox 100") - the "7. pull request" box.
#   * "[Optional:] git branch -d feature" -> split into two lines:
#       "[Optional:] git checkout master"
#       "[Optional:] git branch -D feature"   (note: -d -> -D)
#   * grow the box height to fit the extra line
# ---------------------------------------------------------------
$prShape = $s.Shapes.Item(35)
$prRange = $prShape.TextFrame.TextRange

# Step 1: fix the dash case (-d -> -D) while keeping the paragraph
# count unchanged - this preserves every run's existing formatting.
$prRange.Text = (
    "7. pull request" + $cr +
    "[Use web UI to initiate PR]" + $cr +
    "[Wait for tests on GitHub Actions]" + $cr +
    "[Squash and merge]" + $cr +
    "[Optional:] git branch " + $dash + "D feature" + $cr +
    "[Optional:] git push --delete origin feature"
)

# Step 2: insert the new "git checkout master" paragraph. Changing the
# paragraph count resets run-level formatting (size/italic) across the
# whole text frame, so it is restored explicitly afterwards.
$prRange.Text = (
    "7. pull request" + $cr +
    "[Use web UI to initiate PR]" + $cr +
    "[Wait for tests on GitHub Actions]" + $cr +
    "[Squash and merge]" + $cr +
    "[Optional:] git checkout master" + $cr +
    "[Optional:] git branch " + $dash + "D feature" + $cr +
    "[Optional:] git push --delete origin feature"
)

# Restore the sz=8pt body formatting for every paragraph below the
# sz=14pt title line (chars 17 through the end, char 208).
$prRange.Characters(17, 208 - 17 + 1).Font.Size = 8
# Restore italics on the two "feature" runs.
$prRange.Characters(157, 7).Font.Italic = $true
$prRange.Characters(202, 7).Font.Italic = $true

$prShape.Height = EmuToPt 1046440
